# =====================================================================
# Edit script: updates the "last modified" date fields (12.11.2025 ->
# 14.11.2025 / 11/12/2025 -> 11/14/2025) and a handful of real text
# corrections (produtos -> usuarios route, "Atualiza produto" ->
# "Atualiza usuario"). The remaining hunks in the source diff only
# re-split existing runs (same visible text, PowerPoint proofing
# metadata) and are not reproduced since that metadata is not exposed
# through the PowerPoint object model.
# =====================================================================

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Date placeholder fields
# ---------------------------------------------------------------------

# 1a. Notes master (Czech locale: 12.11.2025 -> 14.11.2025)
$nm = $p.NotesMaster
for ($i = 1; $i -le $nm.Shapes.Count; $i++) {
    $sh = $nm.Shapes.Item($i)
    if ($sh.Name -eq "Date Placeholder 2") {
        $sh.TextFrame.TextRange.Text = "14.11.2025"
    }
}

# 1b. Slide master (English locale: 11/12/2025 -> 11/14/2025)
$sm = $p.SlideMaster
for ($i = 1; $i -le $sm.Shapes.Count; $i++) {
    $sh = $sm.Shapes.Item($i)
    if ($sh.Name -eq "Date Placeholder 3") {
        $sh.TextFrame.TextRange.Text = "11/14/2025"
    }
}

# 1c. All slide layouts (English locale: 11/12/2025 -> 11/14/2025)
for ($li = 1; $li -le $sm.CustomLayouts.Count; $li++) {
    $layout = $sm.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "11/14/2025"
        }
    }
}

Write-Host "Date fields updated"
